# formative-elements.xlsx bug-fix edit:
#   - "suborder" sheet: insert a new row (18) for the "ist" formative
#     element (same tissue/organic-soil-materials meaning as "hist"),
#     shifting the rest of the table down by one row.
#   - Selection / active-sheet bookkeeping updates that go along with the
#     edit: "suborder" becomes the active tab (was "subgroup"), with the
#     cursor left on B14; "subgroup" keeps its own last selection (G14).

$wb = $excel.ActiveWorkbook

$suborder = $wb.Worksheets.Item("suborder")
$subgroup = $wb.Worksheets.Item("subgroup")

# Work on the "suborder" sheet: insert a row above the old row 18 and
# populate it with the new "ist" formative element.
$suborder.Activate() | Out-Null
$suborder.Rows.Item(18).Insert()

$suborder.Range("A18").Value = "ist"
$suborder.Range("B18").Value = "tissue"
$suborder.Range("C18").Value = "presence of organic soil materials"

# Leave the cursor where the author left it after editing this sheet.
$suborder.Range("B14").Select() | Out-Null

# "subgroup" was the active sheet before the edit; make sure it keeps its
# own selection now that it is no longer the active tab.
$subgroup.Activate() | Out-Null
$subgroup.Range("G14").Select() | Out-Null

# "suborder" ends up as the active tab after the edit.
$suborder.Activate() | Out-Null
